# aggiornamento fino a 6 gennaio 2022
# Append new daily rows (465-491) to the data table on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Data rows to append: serial date (col A), col B, col C, col D
$data = @(
    @(44539, 0, 2, 74.93443237167479),
    @(44540, 0, 1, 37.46721618583739),
    @(44541, 0, 1, 37.46721618583739),
    @(44542, 0, 1, 37.46721618583739),
    @(44543, 0, 1, 37.46721618583739),
    @(44544, 0, 1, 37.46721618583739),
    @(44545, 0, 0, 0),
    @(44546, 0, 0, 0),
    @(44547, 0, 0, 0),
    @(44548, 0, 0, 0),
    @(44550, 3, 3, 112.4016485575122),
    @(44551, 1, 4, 149.8688647433496),
    @(44552, 1, 5, 187.3360809291869),
    @(44553, 0, 5, 187.3360809291869),
    @(44554, 0, 5, 187.3360809291869),
    @(44555, 1, 6, 224.8032971150243),
    @(44556, 3, 9, 337.2049456725365),
    @(44557, 8, 14, 524.5410266017235),
    @(44558, 4, 17, 636.9426751592357),
    @(44559, 3, 19, 711.8771075309105),
    @(44560, 5, 24, 899.2131884600974),
    @(44561, 3, 27, 1011.61483701761),
    @(44562, 14, 40, 1498.688647433496),
    @(44563, 4, 41, 1536.155863619333),
    @(44564, 2, 35, 1311.352566504309),
    @(44565, 2, 33, 1236.418134132634),
    @(44566, 5, 35, 1311.352566504309)
)

$startRow = 465

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]

    # Column A carries the same formatting as the rest of the date column
    # (bold, centered, thin box border, custom date/time number format).
    $cellA = $ws.Cells.Item($r, 1)
    $cellA.Value = $row[0]
    $cellA.NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $cellA.Font.Bold = $true
    $cellA.HorizontalAlignment = -4108
    $cellA.VerticalAlignment = -4160
    $cellA.Borders.LineStyle = 1

    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
